$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the date the record was last changed, stored
# as an Excel date serial number. For every data row (2 through 28) this
# value moves forward by one day, from 45181 (2023-09-12) to 45182
# (2023-09-13).
$ws.Range("C2:C28").Value = 45182
